$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("10036957", "HAPPYTOS CHIP HJU140", "PM1MKT", "1", "1", "RT,(E-2B)"),
    @("10003814", "AQUA AIR MINERAL 600", "PM1MKT", "1", "2", "RT,(E-3B)"),
    @("20040383", "NU MILK TEA 330ML",    "PM1MKT", "1", "3", "RT,(E-1B)"),
    @("10036916", "GLICO POCKY.STRAW 45", "PM1MKT", "2", "1", "RT,(E-1B)"),
    @("10008819", "BANGO KECAP MNS 700G", "PM1MKT", "2", "2", "RT,(E-0.5B)"),
    @("20138899", "SNSLK SHP G.BLCK 160", "PM1MKT", "3", "1", "PT,(E-3B)"),
    @("20047217", "PEPSODENT WHITE 225G", "PM1MKT", "3", "2", "PT,(E-3B)"),
    @("20129837", "LARIST SPR.GRD 750ML", "PM1MKT", "4", "1", "PT,(E-1B)"),
    @("20128974", "MONTISS FC TISUE200S", "PM1MKT", "4", "2", "RT")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$c]
    }
}

# The source workbook has one extra trailing row (row 11) that's no longer
# needed after the refresh - remove it so the sheet ends at row 10.
$ws.Rows(11).Delete()
